$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (phone numbers) with the new set
$ws.Range("B1").Value = 81293112671
$ws.Range("B2").Value = 81293112672
$ws.Range("B3").Value = 81293112673
$ws.Range("B4").Value = 81293112674
$ws.Range("B5").Value = 81293112675
$ws.Range("B6").Value = 81293112676

# The "jl. RCM" label column moves from D to E
$ws.Range("D1:D6").Cut($ws.Range("E1:E6"))

# Column B is now widened to best-fit the new (longer) phone numbers
$ws.Columns("B").AutoFit()

# Zoom in to 170%
$excel.ActiveWindow.Zoom = 170

# Selection moves to D1:D6 (now an empty column after the cut)
[void]$ws.Range("D1:D6").Select()
